$d = $word.ActiveDocument

# The document has two "Utilities Req'd: ... IW" lines that are each followed
# by two manual line breaks (<w:br/><w:br/>) before the "See plans..." text.
# The target edit turns the first of those two line breaks into a paragraph
# break instead (splitting what is currently one big paragraph into
# additional paragraphs), while leaving the second line break intact at the
# start of the new paragraph.
#
# "IW" only occurs in the document right before these two spots, so
# searching for "IW" followed by two line breaks (^l^l) reliably finds each
# target location without touching the other (unrelated) double-line-break
# run later in the document.

$maxIterations = 10
for ($i = 0; $i -lt $maxIterations; $i++) {
    $found = $d.Content
    $ok = $found.Find.Execute("IW^l^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        break
    }

    # $found now spans "IW" + the two line-break characters; the first line
    # break is the character right before the last one.
    $breakStart = $found.End - 2
    $breakEnd = $found.End - 1

    # Split the paragraph right after that first line-break character...
    $d.Range($breakStart, $breakEnd).InsertParagraphAfter()
    # ...then remove the (now redundant) line-break character itself, so the
    # paragraph boundary takes its place and the second break starts the new
    # paragraph.
    $d.Range($breakStart, $breakEnd).Delete()
}
